# Fruta / hortaliza, semanal
# Insert two new daily price rows for "Palta" (Hass) at Vega Modelo de Temuco,
# right after the existing 2020-12-07 (serial 44172) block, shifting every
# following row down by two. The two new rows carry date 44491 and the new
# "1a nueva(o)" / "2a nueva(o)" quality grades.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 648 (old row 648 -> 650, old row 649 -> 651, etc.)
$ws.Rows.Item(648).Insert()
$ws.Rows.Item(648).Insert()

# --- New row 648 ---
$ws.Cells.Item(648, 1).Value = 10
$ws.Cells.Item(648, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(648, 3).Value = "La Araucanía"
$ws.Cells.Item(648, 4).Value = 44491
$ws.Cells.Item(648, 5).Value = 9
$ws.Cells.Item(648, 6).Value = "Fruta"
$ws.Cells.Item(648, 7).Value = 100106
$ws.Cells.Item(648, 8).Value = "Oleaginosos"
$ws.Cells.Item(648, 9).Value = 100106002
$ws.Cells.Item(648, 10).Value = "Palta"
$ws.Cells.Item(648, 11).Value = "Hass"
$ws.Cells.Item(648, 12).Value = "1a nueva(o)"
$ws.Cells.Item(648, 13).Value = 300
$ws.Cells.Item(648, 14).Value = 2800
$ws.Cells.Item(648, 15).Value = 3000
$ws.Cells.Item(648, 16).Value = 2900
$ws.Cells.Item(648, 17).Value = "$/kilo (en bandeja de 18 kilos)"
$ws.Cells.Item(648, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(648, 19).Value = 2900
$ws.Cells.Item(648, 20).Value = 1

# --- New row 649 ---
$ws.Cells.Item(649, 1).Value = 10
$ws.Cells.Item(649, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(649, 3).Value = "La Araucanía"
$ws.Cells.Item(649, 4).Value = 44491
$ws.Cells.Item(649, 5).Value = 9
$ws.Cells.Item(649, 6).Value = "Fruta"
$ws.Cells.Item(649, 7).Value = 100106
$ws.Cells.Item(649, 8).Value = "Oleaginosos"
$ws.Cells.Item(649, 9).Value = 100106002
$ws.Cells.Item(649, 10).Value = "Palta"
$ws.Cells.Item(649, 11).Value = "Hass"
$ws.Cells.Item(649, 12).Value = "2a nueva(o)"
$ws.Cells.Item(649, 13).Value = 100
$ws.Cells.Item(649, 14).Value = 2500
$ws.Cells.Item(649, 15).Value = 2500
$ws.Cells.Item(649, 16).Value = 2500
$ws.Cells.Item(649, 17).Value = "$/kilo (en bandeja de 18 kilos)"
$ws.Cells.Item(649, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(649, 19).Value = 2500
$ws.Cells.Item(649, 20).Value = 1
